# Added raid levels to spreadsheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "TemplateSheet"
$ws.Name = "TemplateSheet"

# --- Column A / D labels -----------------------------------------------
# (Set in this order so the shared-string table is built with
#  "Number of Drives" before "RAID", matching the source workbook.)
$ws.Range("A2").Value = "Number of Drives"

$ws.Range("A1").Value = "RAID"
$ws.Range("A1").Font.Bold = $true

$ws.Range("A3").Value = "Single Drive Capacity"
$ws.Range("A4").Value = "# Parity Drives"

$ws.Range("D1").Value = "Index calculations:"
$ws.Range("D2").Value = "Disk block size:"
$ws.Range("D3").Value = "Record size:"
$ws.Range("D4").Value = "Record size:"
$ws.Range("D5").Value = "Data File Size: "

# --- RAID level numbers (column A, rows 5-8) ----------------------------
$ws.Range("A5").Value = 0
$ws.Range("A6").Value = 1
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6

# --- Inputs and usable-capacity formulas (column B) ---------------------
$ws.Range("B2").Value = 8      # Number of Drives
$ws.Range("B3").Value = 8000   # Single Drive Capacity
$ws.Range("B4").Value = 0      # # Parity Drives

$ws.Range("B5").Formula = "=B2*B3"
$ws.Range("B6").Formula = "=CEILING.MATH(B2/2)*B3"
$ws.Range("B7").Formula = "=(B2-1)*B3"
$ws.Range("B8").Formula = "=(B2-B4)*B3"

# --- Column widths --------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 21.90625
$ws.Columns.Item(4).ColumnWidth = 17.26953125

# --- Final selection, matching the saved view state ----------------------
$ws.Range("D5").Select() | Out-Null
